# edit.ps1 - Apply the "2-x." -> "1-x." title renumbering across the
# module/section title slides, and expand the title-slide subtitle run
# "Introduction to RNA sequencing (lecture)" into
# "Introduction to RNA sequencing (tutorial)".
#
# All of the "2-..." titles live in Shape 1 (the title/rectangle
# placeholder) of their slide, as a single run starting "2-<roman>. ...".
# We split that single run into three runs: "1", "-", "<roman>. ...",
# by rewriting the first two characters in place (this is exactly the
# run-split pattern PowerPoint itself performs when a user edits text
# in place).

$p = $ppt.ActivePresentation

# Slide index -> shape index holding the "2-...” title text.
$slides = @(10, 11, 12, 6, 7, 8, 9)

foreach ($slideIdx in $slides) {
    $slide = $p.Slides.Item($slideIdx)
    $shape = $slide.Shapes.Item(1)
    $tr = $shape.TextFrame.TextRange

    # First character: "2" -> "1"
    $first = $tr.Characters(1, 1)
    $first.Text = "1"

    # Second character stays "-" but gets split into its own run.
    $second = $tr.Characters(2, 1)
    $second.Text = "-"
}

# Slide 2 title-slide subtitle: "...Introduction to RNA sequencing (lecture)"
# -> "...Introduction to RNA sequencing (tutorial)", with "(", "tutorial",
# ")" becoming their own runs.
$slide2 = $p.Slides.Item(2)
$shape2 = $slide2.Shapes.Item(4)
$tr2 = $shape2.TextFrame.TextRange

$fullText = $tr2.Text
$openParenPos = $fullText.IndexOf("(lecture)") + 1

$openParen = $tr2.Characters($openParenPos, 1)
$openParen.Text = "("

$word = $tr2.Characters($openParenPos + 1, 7)
$word.Text = "tutorial"

$closeParen = $tr2.Characters($openParenPos + 9, 1)
$closeParen.Text = ")"
